# Apply updated TPM-derived values to LR-pairs sheet (Vim-Cd44)
function Set-RowRange {
    param($ws, $range, $values)
    $arr = New-Object "object[,]" 1,$values.Count
    for ($i = 0; $i -lt $values.Count; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($range).Value = $arr
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-RowRange $ws "G2:J2" @(1545.224243, 4635.672729, 0.3401139553078116, 0.3401139553078116)
Set-RowRange $ws "M2:T2" @(19.21315233333334, 57.63945700000001, 0.04451179209991234, 0.04451179209991233, 29688.62876991869, 267197.6589292682, 0.01513908166894018, 0.01513908166894018)

Set-RowRange $ws "G3:J3" @(1545.224243, 4635.672729, 0.3401139553078116, 0.3401139553078116)
Set-RowRange $ws "O3:T3" @(0.2141755495962477, 0.2141755495962477, 142851.5475019221, 1285663.927517299, 0.07284409330340419, 0.07284409330340419)

Set-RowRange $ws "G4:J4" @(1545.224243, 4635.672729, 0.3401139553078116, 0.3401139553078116)
Set-RowRange $ws "M4:T4" @(166.8580016666666, 500.5740049999999, 0.3865658561145097, 0.3865658561145097, 257833.0293138677, 2320497.263824809, 0.1314764423100563, 0.1314764423100563)

Set-RowRange $ws "G5:J5" @(1545.224243, 4635.672729, 0.3401139553078116, 0.3401139553078116)
Set-RowRange $ws "M5:T5" @(41.09915599999999, 123.297468, 0.09521587377309249, 0.09521587377309249, 63507.41221803889, 571566.7099623501, 0.03238424743705581, 0.03238424743705581)

Set-RowRange $ws "G6:J6" @(1545.224243, 4635.672729, 0.3401139553078116, 0.3401139553078116)
Set-RowRange $ws "M6:T6" @(112.0244103333333, 336.073231, 0.2595309284162377, 0.2595309284162377, 173102.8346548464, 1557925.511893617, 0.08827009058835511, 0.08827009058835512)

Set-RowRange $ws "I7:J7" @(0.2751988005753909, 0.2751988005753909)
Set-RowRange $ws "M7:T7" @(19.21315233333334, 57.63945700000001, 0.04451179209991234, 0.04451179209991233, 24022.16933678997, 216199.5240311097, 0.01224959179735703, 0.01224959179735704)

Set-RowRange $ws "I8:J8" @(0.2751988005753909, 0.2751988005753909)
Set-RowRange $ws "O8:P8" @(0.2141755495962477, 0.2141755495962477)
Set-RowRange $ws "S8:T8" @(0.05894085436146252, 0.05894085436146252)

Set-RowRange $ws "I9:J9" @(0.2751988005753909, 0.2751988005753909)
Set-RowRange $ws "M9:T9" @(166.8580016666666, 500.5740049999999, 0.3865658561145097, 0.3865658561145097, 208622.2553016963, 1877600.297715267, 0.1063824599461122, 0.1063824599461122)

Set-RowRange $ws "I10:J10" @(0.2751988005753909, 0.2751988005753909)
Set-RowRange $ws "M10:T10" @(41.09915599999999, 123.297468, 0.09521587377309249, 0.09521587377309249, 51386.19982303861, 462475.7984073475, 0.02620329425809287, 0.02620329425809288)

Set-RowRange $ws "I11:J11" @(0.2751988005753909, 0.2751988005753909)
Set-RowRange $ws "M11:T11" @(112.0244103333333, 336.073231, 0.2595309284162377, 0.2595309284162377, 140063.9160192666, 1260575.2441734, 0.07142260021236624, 0.07142260021236625)

Set-RowRange $ws "G12:J12" @(618.110189, 1854.330567, 0.1360500925022369, 0.1360500925022369)
Set-RowRange $ws "M12:T12" @(19.21315233333334, 57.63945700000001, 0.04451179209991234, 0.04451179209991233, 11875.84522004246, 106882.6069803821, 0.00605583343263341, 0.00605583343263341)

Set-RowRange $ws "G13:J13" @(618.110189, 1854.330567, 0.1360500925022369, 0.1360500925022369)
Set-RowRange $ws "O13:T13" @(0.2141755495962477, 0.2141755495962477, 57142.51340888103, 514282.6206799292, 0.02913860333428693, 0.02913860333428693)

Set-RowRange $ws "G14:J14" @(618.110189, 1854.330567, 0.1360500925022369, 0.1360500925022369)
Set-RowRange $ws "M14:T14" @(166.8580016666666, 500.5740049999999, 0.3865658561145097, 0.3865658561145097, 103136.6309463456, 928229.6785171108, 0.05259232048258543, 0.05259232048258545)

Set-RowRange $ws "G15:J15" @(618.110189, 1854.330567, 0.1360500925022369, 0.1360500925022369)
Set-RowRange $ws "M15:T15" @(41.09915599999999, 123.297468, 0.09521587377309249, 0.09521587377309249, 25403.80708290048, 228634.2637461043, 0.01295412843451054, 0.01295412843451055)

Set-RowRange $ws "G16:J16" @(618.110189, 1854.330567, 0.1360500925022369, 0.1360500925022369)
Set-RowRange $ws "M16:T16" @(112.0244103333333, 336.073231, 0.2595309284162377, 0.2595309284162377, 69243.42944375022, 623190.864993752, 0.03530920681822056, 0.03530920681822056)

Set-RowRange $ws "G17:J17" @(528.755229, 1586.265687, 0.1163824817915944, 0.1163824817915944)
Set-RowRange $ws "M17:T17" @(19.21315233333334, 57.63945700000001, 0.04451179209991234, 0.04451179209991233, 10159.05476182355, 91431.49285641198, 0.005180392833579281, 0.005180392833579282)

Set-RowRange $ws "G18:J18" @(528.755229, 1586.265687, 0.1163824817915944, 0.1163824817915944)
Set-RowRange $ws "O18:T18" @(0.2141755495962477, 0.2141755495962477, 48881.90374604626, 439937.1337144163, 0.02492628200109002, 0.02492628200109002)

Set-RowRange $ws "G19:J19" @(528.755229, 1586.265687, 0.1163824817915944, 0.1163824817915944)
Set-RowRange $ws "M19:T19" @(166.8580016666666, 500.5740049999999, 0.3865658561145097, 0.3865658561145097, 88227.0408817407, 794043.3679356663, 0.04498949371049901, 0.04498949371049902)

Set-RowRange $ws "G20:J20" @(528.755229, 1586.265687, 0.1163824817915944, 0.1163824817915944)
Set-RowRange $ws "M20:T20" @(41.09915599999999, 123.297468, 0.09521587377309249, 0.09521587377309249, 21731.39364248672, 195582.5427823805, 0.01108145969566768, 0.01108145969566769)

Set-RowRange $ws "G21:J21" @(528.755229, 1586.265687, 0.1163824817915944, 0.1163824817915944)
Set-RowRange $ws "M21:T21" @(112.0244103333333, 336.073231, 0.2595309284162377, 0.2595309284162377, 59233.49273939163, 533101.4346545248, 0.03020485355075836, 0.03020485355075837)

Set-RowRange $ws "G22:J22" @(600.866618, 1802.599854, 0.1322546698229662, 0.1322546698229662)
Set-RowRange $ws "M22:T22" @(19.21315233333334, 57.63945700000001, 0.04451179209991234, 0.04451179209991233, 11544.54186364881, 103900.8767728393, 0.005886892367402422, 0.005886892367402421)

Set-RowRange $ws "G23:J23" @(600.866618, 1802.599854, 0.1322546698229662, 0.1322546698229662)
Set-RowRange $ws "O23:T23" @(0.2141755495962477, 0.2141755495962477, 55548.39474748409, 499935.5527273567, 0.02832571659600407, 0.02832571659600407)

Set-RowRange $ws "G24:J24" @(600.866618, 1802.599854, 0.1322546698229662, 0.1322546698229662)
Set-RowRange $ws "M24:T24" @(166.8580016666666, 500.5740049999999, 0.3865658561145097, 0.3865658561145097, 100259.4031476884, 902334.6283291952, 0.05112513966525675, 0.05112513966525675)

Set-RowRange $ws "G25:J25" @(600.866618, 1802.599854, 0.1322546698229662, 0.1322546698229662)
Set-RowRange $ws "M25:T25" @(41.09915599999999, 123.297468, 0.09521587377309249, 0.09521587377309249, 24695.11086837441, 222255.9978153697, 0.01259274394776557, 0.01259274394776557)

Set-RowRange $ws "G26:J26" @(600.866618, 1802.599854, 0.1322546698229662, 0.1322546698229662)
Set-RowRange $ws "M26:T26" @(112.0244103333333, 336.073231, 0.2595309284162377, 0.2595309284162377, 173102.8346548464, 1557925.511893617, 0.0343241772465374, 0.0343241772465374)
